$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "33.949.77"
Set-TextValue "E2" "  -0.35%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.785.71"
Set-TextValue "E3" "  -0.23%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.11%  "

# Row 5 - BNB
Set-TextValue "D5" "226.57"
Set-TextValue "E5" "  +2.15%  "

# Row 6 - XRP
Set-TextValue "E6" "  -1.48%  "

# Row 7 - USDC
Set-TextValue "E7" "  +0.08%  "

# Row 8 - Solana
Set-TextValue "D8" "32.17"
Set-TextValue "E8" "  -1.16%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.294"
Set-TextValue "E9" "  +3.60%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0681"
Set-TextValue "E10" "  -4.02%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0939"
Set-TextValue "E11" "  +1.26%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "2.040.36"
Set-TextValue "E12" "  -0.29%  "

# Row 13 - Chainlink
Set-TextValue "E13" "  +2.18%  "

# Row 14 - WrappedEther
Set-TextValue "D14" "1.777.84"
Set-TextValue "E14" "  -0.73%  "

# Row 15 - WrappedBTC
Set-TextValue "D15" "33.915.28"
Set-TextValue "E15" "  -0.44%  "

# Row 16 - Polygon
Set-TextValue "D16" "0.617"
Set-TextValue "E16" "  -1.56%  "

# Row 17 - Polkadot
Set-TextValue "E17" "  +0.12%  "

# Row 18 - Litecoin
Set-TextValue "D18" "67.64"
Set-TextValue "E18" "  -0.46%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "242.07"

# Row 20 - ShibaInu
Set-TextValue "D20" "0.0₃0771"
Set-TextValue "E20" "  -1.24%  "

# Row 21 - Dai
Set-TextValue "E21" "  +0.16%  "

# Row 22 - Avalanche
Set-TextValue "D22" "10.67"
Set-TextValue "E22" "  -1.18%  "

# Row 23 - Uniswap
Set-TextValue "E23" "  -0.33%  "

# Row 24 - Toncoin
Set-TextValue "E24" "  -2.35%  "

# Row 25 - Monero
Set-TextValue "D25" "161.80"
Set-TextValue "E25" "  +2.42%  "

# Row 26 - Cosmos
Set-TextValue "E26" "  +1.21%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "16.16"
Set-TextValue "E27" "  -1.32%  "

# Row 28 - Stellar
Set-TextValue "D28" "0.113"

# Row 29 - BinanceUSD
Set-TextValue "E29" "  +0.22%  "

# Row 30 - PancakeSwap
Set-TextValue "E30" "  +2.55%  "

# Row 31 - Hedera
Set-TextValue "E31" "  -1.18%  "

# Row 32 - Filecoin
Set-TextValue "D32" "3.63"
Set-TextValue "E32" "  -1.20%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue "E33" "  +1.98%  "

# Row 34 - LidoDAOToken
Set-TextValue "E34" "  +1.09%  "

# Row 35 - Maker
Set-TextValue "D35" "1.397.44"
Set-TextValue "E35" "  +0.10%  "

# Row 36 - ImmutableX
Set-TextValue "D36" "0.648"
Set-TextValue "E36" "  +1.24%  "

# Row 37 - TrustWalletToken
Set-TextValue "E37" "  -1.23%  "

# Row 38 - was VeChain, now RenderToken
Set-TextValue "B38" "RenderToken"
Set-TextValue "C38" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D38" "2.35"
Set-TextValue "E38" "  +8.72%  "

# Row 39 - was RenderToken, now VeChain
Set-TextValue "B39" "VeChain"
Set-TextValue "C39" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D39" "0.0187"
Set-TextValue "E39" "  +1.19%  "

# Row 40 - Aave
Set-TextValue "D40" "79.87"
Set-TextValue "E40" "  +0.29%  "

# Row 41 - HuobiToken
Set-TextValue "E41" "  +0.46%  "

# Row 42 - ARBITRUM
Set-TextValue "D42" "0.920"
Set-TextValue "E42" "  -0.30%  "

# Row 43 - InjectiveProtocol
Set-TextValue "D43" "13.65"
Set-TextValue "E43" "  +13.62%  "

# Row 44 - MXToken
Set-TextValue "E44" "  -1.61%  "

# Row 45 - BabyDogeCoin
Set-TextValue "E45" "  +8.75%  "

# Row 46 - Kaspa
Set-TextValue "D46" "0.0509"
Set-TextValue "E46" "  +3.21%  "

# Row 47 - WEMIXToken
Set-TextValue "E47" "  +2.79%  "

# Row 48 - FraxShare
Set-TextValue "D48" "5.93"
Set-TextValue "E48" "  +0.76%  "

# Row 49 - Quant
Set-TextValue "D49" "107.60"
Set-TextValue "E49" "  +0.28%  "

# Row 50 - RocketPoolETH
Set-TextValue "D50" "1.942.50"
Set-TextValue "E50" "  -0.48%  "

# Row 51 - PaxDollar
Set-TextValue "E51" "  +0.17%  "
